# Applies the 2024-02-24 23:13 data update to "Spain Primera Liga" sheet.
# Two kinds of edits:
#  1) Several row-pairs had their data (columns B, F..AC) swapped between the
#     two rows (column A id, and columns C/D/E stay put on their own row).
#  2) A handful of single cells in rows 485/488/493 were corrected in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell($r1, $r2, $c) {
    $a = $ws.Cells.Item($r1, $c).Value2
    $b = $ws.Cells.Item($r2, $c).Value2
    $ws.Cells.Item($r1, $c).Value2 = $b
    $ws.Cells.Item($r2, $c).Value2 = $a
}

function Swap-Rows($rowA, $rowB) {
    # Columns that participate in the swap: B, and F through AC
    # (C, D, E are left untouched on each row).
    $cols = @(2) + @(6..29)
    foreach ($c in $cols) {
        Swap-Cell $rowA $rowB $c
    }
}

# Row pairs whose content got swapped.
Swap-Rows 388 389
Swap-Rows 404 405
Swap-Rows 406 407
Swap-Rows 408 409

# Direct single-cell corrections.
$ws.Range("U485").Value2 = 2
$ws.Range("V485").Value2 = 1.9

$ws.Range("P488").Value2 = 8.5
$ws.Range("Q488").Value2 = -1.5
$ws.Range("R488").Value2 = 1.84
$ws.Range("S488").Value2 = 2.06

$ws.Range("R493").Value2 = 2.06
$ws.Range("S493").Value2 = 1.84
